$wb = $excel.ActiveWorkbook

# Sheet "展览" (index 1 / sheetId 1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F7").Value = 599
$ws1.Range("F9").Value = 8858
$ws1.Range("F10").Value = 817
$ws1.Range("F11").Value = 334
$ws1.Range("F12").Value = 1149
$ws1.Range("F13").Value = 1004
$ws1.Range("F14").Value = 118
$ws1.Range("F21").Value = 1076

# Sheet "全部类型" (index 4 / sheetId 4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value = 599
$ws4.Range("F11").Value = 8858
$ws4.Range("F12").Value = 817
$ws4.Range("F13").Value = 334
$ws4.Range("F14").Value = 1149
$ws4.Range("F15").Value = 1004
$ws4.Range("F16").Value = 118
$ws4.Range("F23").Value = 1076
